$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (on what were originally slides 14, 15 and 16) get their
#    table style switched from the deck's custom "Table_0" style
#    ({D4014D65-F22D-4818-AA95-E53A00D17609}) to the built-in "Medium Style 2
#    - Accent 1" style ({8561FC4F-442F-4506-A984-9B81CC9B8A49}). Each table is
#    the first shape (a graphicFrame) on its slide.
# ---------------------------------------------------------------------------
$newTableStyleId = "{8561FC4F-442F-4506-A984-9B81CC9B8A49}"

$tableSlideIndexes = @(14, 15, 16)
foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newTableStyleId)
}

# ---------------------------------------------------------------------------
# 2) The slide master's theme ("Integral" / "Red Violet" colour scheme) is
#    replaced with the stock Office theme's colour scheme. The font scheme
#    and format scheme are unchanged (they were already identical between the
#    two themes embedded in this deck), so only the twelve theme colours
#    need to be rewritten. The colour-scheme slot order matches the
#    <a:clrScheme> child order: dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink. RGB() is unavailable in this host, so colours are supplied as
#    the equivalent OLE BGR integers.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$master.ColorScheme.Colors(1).RGB = 0          # dk1      000000
$master.ColorScheme.Colors(2).RGB = 16777215   # lt1      FFFFFF
$master.ColorScheme.Colors(3).RGB = 6968388    # dk2      44546A
$master.ColorScheme.Colors(4).RGB = 15132391   # lt2      E7E6E6
$master.ColorScheme.Colors(5).RGB = 13998939   # accent1  5B9BD5
$master.ColorScheme.Colors(6).RGB = 3243501    # accent2  ED7D31
$master.ColorScheme.Colors(7).RGB = 10855845   # accent3  A5A5A5
$master.ColorScheme.Colors(8).RGB = 49407      # accent4  FFC000
$master.ColorScheme.Colors(9).RGB = 12874308   # accent5  4472C4
$master.ColorScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$master.ColorScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$master.ColorScheme.Colors(12).RGB = 7491477   # folHlink 954F72
